$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.682.69"
$ws.Range("E2").Value = "  +5.20%  "
$ws.Range("D3").Value = "2.225.43"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "61.06"
$ws.Range("E7").Value = "  -3.18%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  +3.18%  "
$ws.Range("D10").Value = "58.18"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "0.0881"
$ws.Range("E11").Value = "  +4.30%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "2.557.89"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").Value = "15.57"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "21.48"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "0.793"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "5.53"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "2.230.10"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").Value = "41.555.88"
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("D20").Value = "72.53"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  +5.35%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "248.97"
$ws.Range("E23").Value = "  +8.49%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  -5.61%  "
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "167.75"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "0.141"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").Value = "19.89"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").Value = "2.58"
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  +7.57%  "
$ws.Range("D35").Value = "4.62"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").Value = "0.0623"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "6.55"
$ws.Range("E37").Value = "  -4.99%  "
$ws.Range("D38").Value = "3.67"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "0.000238"
$ws.Range("E41").Value = "  +27.70%  "
$ws.Range("D42").Value = "4.85"
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("D44").Value = "8.68"
$ws.Range("E44").Value = "  +11.75%  "
$ws.Range("E45").Value = "  +7.10%  "
$ws.Range("D46").Value = "98.93"
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.464.88"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").Value = "16.33"
$ws.Range("E49").Value = "  -7.31%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  -1.13%  "
